{"js": "// Submit quiz7-8, Associated Press-style headlines:\n// \"Design, Implementation & Software Testing\" -> \"Design, implementation and software testing\"\nconst body = context.document.body;\n\n// Narrow the edit to the exact run of text that holds the homework title so\n// the leading \"Homework 7-8: \" (hyperlink + colon run) is left untouched.\nconst results = body.search(\"Design, Implementation & Software Testing\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target heading text \"Design, Implementation & Software Testing\" was not found.');\n}\n\nconst target = results.items[0];\n// Replace in place so the run keeps its existing bold/size/theme-font formatting.\ntarget.insertText(\"Design, implementation and software testing\", \"Replace\");\nawait context.sync();\n", "ps1": "# Submit quiz7-8, Associated Press-style headlines:\n# \"Design, Implementation & Software Testing\" -> \"Design, implementation and software testing\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Design, Implementation & Software Testing\"\n$find.Replacement.Text = \"Design, implementation and software testing\"\n\n# FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n# MatchAllWordForms, Forward, Wrap (0 = wdFindStop), Format, ReplaceWith,\n# Replace (1 = wdReplaceOne) \u2014 only the single heading occurrence changes.\n$find.Execute(\n    [ref]$find.Text,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]0,\n    [ref]$false,\n    [ref]$find.Replacement.Text,\n    [ref]1\n) | Out-Null\n"}
